$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 19 / A19: centered date style (new style #5) ---
$ws.Cells.Item(19,1).Value = 44076
$ws.Cells.Item(19,1).HorizontalAlignment = -4108
$ws.Cells.Item(19,1).NumberFormat = "mm-dd-yy"

# --- Block starting row 19: Scrum ---
$ws.Cells.Item(19,2).Value = "Scrum"
$ws.Cells.Item(19,3).Value = 1
$ws.Cells.Item(4,3).Copy()
$ws.Cells.Item(19,3).PasteSpecial(-4122)
$ws.Cells.Item(19,4).Value = "https://youtu.be/q0QL9jbSx4c"
$ws.Hyperlinks.Add($ws.Cells.Item(19,4), "https://youtu.be/q0QL9jbSx4c") | Out-Null
$ws.Cells.Item(20,3).Value = 2
$ws.Cells.Item(4,3).Copy()
$ws.Cells.Item(20,3).PasteSpecial(-4122)
$ws.Cells.Item(20,4).Value = "https://youtu.be/V85VSxP9Hz4"
$ws.Hyperlinks.Add($ws.Cells.Item(20,4), "https://youtu.be/V85VSxP9Hz4") | Out-Null
$ws.Cells.Item(21,3).Value = 3
$ws.Cells.Item(4,3).Copy()
$ws.Cells.Item(21,3).PasteSpecial(-4122)
$ws.Cells.Item(21,4).Value = "https://youtu.be/KXtjGzz-AwM"
$ws.Hyperlinks.Add($ws.Cells.Item(21,4), "https://youtu.be/KXtjGzz-AwM") | Out-Null

# --- Block starting row 23: Sprint ---
$ws.Cells.Item(23,2).Value = "Sprint"
$ws.Cells.Item(23,1).Value = 44083
$ws.Cells.Item(23,1).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(23,3).Value = 1
$ws.Cells.Item(4,3).Copy()
$ws.Cells.Item(23,3).PasteSpecial(-4122)
$ws.Cells.Item(23,4).Value = "https://youtu.be/sjpEWQ61R5o"
$ws.Hyperlinks.Add($ws.Cells.Item(23,4), "https://youtu.be/sjpEWQ61R5o") | Out-Null
$ws.Cells.Item(24,3).Value = 2
$ws.Cells.Item(4,3).Copy()
$ws.Cells.Item(24,3).PasteSpecial(-4122)
$ws.Cells.Item(24,4).Value = "https://youtu.be/C5FVpLiZA24"
$ws.Hyperlinks.Add($ws.Cells.Item(24,4), "https://youtu.be/C5FVpLiZA24") | Out-Null

# --- Block starting row 26: Testing de Caja Negra ---
$ws.Cells.Item(26,2).Value = "Testing de Caja Negra"
$ws.Cells.Item(26,1).Value = 44090
$ws.Cells.Item(23,1).Copy()
$ws.Cells.Item(26,1).PasteSpecial(-4122)
$ws.Cells.Item(26,3).Value = 1
$ws.Cells.Item(4,3).Copy()
$ws.Cells.Item(26,3).PasteSpecial(-4122)
$ws.Cells.Item(26,4).Value = "https://youtu.be/OUY0N9cuz18"
$ws.Hyperlinks.Add($ws.Cells.Item(26,4), "https://youtu.be/OUY0N9cuz18") | Out-Null
$ws.Cells.Item(27,3).Value = 2
$ws.Cells.Item(4,3).Copy()
$ws.Cells.Item(27,3).PasteSpecial(-4122)
$ws.Cells.Item(27,4).Value = "https://youtu.be/QJkThF0MpDs"
$ws.Hyperlinks.Add($ws.Cells.Item(27,4), "https://youtu.be/QJkThF0MpDs") | Out-Null

# --- Block starting row 29: Testing Caja Blanca ---
$ws.Cells.Item(29,2).Value = "Testing Caja Blanca"
$ws.Cells.Item(29,1).Value = 44097
$ws.Cells.Item(23,1).Copy()
$ws.Cells.Item(29,1).PasteSpecial(-4122)
$ws.Cells.Item(29,3).Value = 1
$ws.Cells.Item(4,3).Copy()
$ws.Cells.Item(29,3).PasteSpecial(-4122)
$ws.Cells.Item(29,4).Value = "https://youtu.be/5kfUgtONLE0"
$ws.Hyperlinks.Add($ws.Cells.Item(29,4), "https://youtu.be/5kfUgtONLE0") | Out-Null
$ws.Cells.Item(30,3).Value = 2
$ws.Cells.Item(4,3).Copy()
$ws.Cells.Item(30,3).PasteSpecial(-4122)
$ws.Cells.Item(30,4).Value = "https://youtu.be/o0IJSCVb7t8"
$ws.Hyperlinks.Add($ws.Cells.Item(30,4), "https://youtu.be/o0IJSCVb7t8") | Out-Null

# --- Block starting row 32: Ejecución de Casos de Prueba ---
$ws.Cells.Item(32,2).Value = "Ejecución de Casos de Prueba"
$ws.Cells.Item(32,1).Value = 44111
$ws.Cells.Item(23,1).Copy()
$ws.Cells.Item(32,1).PasteSpecial(-4122)
$ws.Cells.Item(32,3).Value = 1
$ws.Cells.Item(4,3).Copy()
$ws.Cells.Item(32,3).PasteSpecial(-4122)
$ws.Cells.Item(32,4).Value = "https://youtu.be/4sxrbciyBZY"
$ws.Hyperlinks.Add($ws.Cells.Item(32,4), "https://youtu.be/4sxrbciyBZY") | Out-Null

# --- Block starting row 34: Repaso para el parcial ---
$ws.Cells.Item(34,2).Value = "Repaso para el parcial"
$ws.Cells.Item(34,1).Value = 37544
$ws.Cells.Item(23,1).Copy()
$ws.Cells.Item(34,1).PasteSpecial(-4122)
$ws.Cells.Item(34,3).Value = 1
$ws.Cells.Item(4,3).Copy()
$ws.Cells.Item(34,3).PasteSpecial(-4122)
$ws.Cells.Item(34,4).Value = "https://youtu.be/jK4MykdHH40"
$ws.Hyperlinks.Add($ws.Cells.Item(34,4), "https://youtu.be/jK4MykdHH40") | Out-Null

# --- Block starting row 36: Practico 13 Design Thinking ---
$ws.Cells.Item(36,2).Value = "Practico 13 Design Thinking"
$ws.Cells.Item(36,1).Value = 44125
$ws.Cells.Item(23,1).Copy()
$ws.Cells.Item(36,1).PasteSpecial(-4122)
$ws.Cells.Item(36,3).Value = 1
$ws.Cells.Item(36,3).HorizontalAlignment = -4108
$ws.Cells.Item(36,4).Value = "https://youtu.be/ZLYnX0E4Uf4"
$ws.Hyperlinks.Add($ws.Cells.Item(36,4), "https://youtu.be/ZLYnX0E4Uf4") | Out-Null

$ws.Application.CutCopyMode = $false

# --- Final selection / view state ---
$ws.Range("F41").Select()